# Apply the cryptocurrency price/volume update described in the commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.580.07'
$ws.Range("E2").Value = '  +2.72%  '
$ws.Range("D3").Value = '1.670.44'
$ws.Range("E3").Value = '  +2.29%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9989'
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.85'
$ws.Range("E5").Value = '  +0.98%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9999'
$ws.Range("E6").Value = '  -0.06%  '
$ws.Range("E7").Value = '  +1.28%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06172'
$ws.Range("E9").Value = '  +2.46%  '
$ws.Range("D10").Value = '1.671.83'
$ws.Range("E10").Value = '  +2.33%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06992'
$ws.Range("E11").Value = '  -0.33%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.81'
$ws.Range("E12").Value = '  +0.69%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.5890'
$ws.Range("E13").Value = '  -3.93%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.383'
$ws.Range("E14").Value = '  +0.81%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '75.34'
$ws.Range("E15").Value = '  +3.93%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.0000'
$ws.Range("E16").Value = '  -0.07%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9998'
$ws.Range("E17").Value = '  +0.15%  '
$ws.Range("D18").Value = '25.568.80'
$ws.Range("E18").Value = '  +2.61%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000006749'
$ws.Range("E19").Value = '  +3.42%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.43'
$ws.Range("E20").Value = '  +3.49%  '
$ws.Range("D21").Value = '1.888.41'
$ws.Range("E21").Value = '  +2.36%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.451'
$ws.Range("E22").Value = '  +2.00%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.779'
$ws.Range("E23").Value = '  +2.59%  '
$ws.Range("E24").Value = '  +0.18%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '136.34'
$ws.Range("E25").Value = '  +1.45%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.04'
$ws.Range("E26").Value = '  +2.13%  '
$ws.Range("E27").Value = '  +1.35%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.724'
$ws.Range("E28").Value = '  +5.22%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '104.68'
$ws.Range("E29").Value = '  +1.93%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.983'
$ws.Range("E30").Value = '  +6.25%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.07870'
$ws.Range("E31").Value = '  +2.01%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.631'
$ws.Range("E32").Value = '  +2.58%  '
$ws.Range("B33").Value = 'Frax'
$ws.Range("C33").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9990'
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04291'
$ws.Range("E34").Value = '  +0.26%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.621'
$ws.Range("E35").Value = '  +1.03%  '
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9576'
$ws.Range("E36").Value = '  +4.27%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6061'
$ws.Range("E37").Value = '  +5.07%  '
$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.590'
$ws.Range("E38").Value = '  +1.26%  '
$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.9061'
$ws.Range("E39").Value = '  +9.73%  '
$ws.Range("B40").Value = 'PaxDollar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9998'
$ws.Range("E40").Value = '  +0.21%  '
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.860'
$ws.Range("E41").Value = '  +4.06%  '
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.01478'
$ws.Range("E42").Value = '  -4.12%  '
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '97.72'
$ws.Range("E43").Value = '  +0.87%  '
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3761'
$ws.Range("E44").Value = '  +1.93%  '
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.885'
$ws.Range("E45").Value = '  +3.63%  '
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1117'
$ws.Range("E46").Value = '  +1.92%  '
$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.220'
$ws.Range("E47").Value = '  +2.90%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05265'
$ws.Range("E48").Value = '  +1.08%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '29.93'
$ws.Range("E49").Value = '  +1.68%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.451'
$ws.Range("E50").Value = '  +4.58%  '
$ws.Range("B51").Value = 'TrueUSD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.001'
$ws.Range("E51").Value = '  +0.21%  '
